$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("DBD")

$ws.Range("D14").Value = "TIMESTAMP"
$ws.Range("E14").Clear()

$ws.Range("D16").Value = "TIMESTAMP"
$ws.Range("E16").Clear()
